# Refresh the live crypto price/volume figures (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.218.08"
$ws.Range("E2").Value = "  +6.58%  "
$ws.Range("D3").Value = "3.019.23"
$ws.Range("E3").Value = "  +3.70%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'585.54"
$ws.Range("E5").Value = "  +2.83%  "
$ws.Range("D6").Value = "'163.21"
$ws.Range("E6").Value = "  +13.30%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.015.51"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").Value = "'6.72"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").Value = "  +5.62%  "
$ws.Range("E12").Value = "  +6.33%  "
$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  +8.61%  "
$ws.Range("D14").Value = "'34.85"
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "66.147.03"
$ws.Range("E16").Value = "  +6.53%  "
$ws.Range("D17").Value = "3.517.53"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("E18").Value = "  +6.91%  "
$ws.Range("D19").Value = "3.021.09"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D20").Value = "'459.19"
$ws.Range("E20").Value = "  +6.54%  "
$ws.Range("D21").Value = "'13.97"
$ws.Range("E21").Value = "  +6.96%  "
$ws.Range("E22").Value = "  +5.85%  "
$ws.Range("E23").Value = "  +8.00%  "
$ws.Range("D24").Value = "'82.48"
$ws.Range("E24").Value = "  +4.50%  "
$ws.Range("D25").Value = "'2.32"
$ws.Range("E25").Value = "  +14.68%  "
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'8.12"
$ws.Range("E29").Value = "  +16.66%  "
$ws.Range("D30").Value = "'2.39"
$ws.Range("E30").Value = "  +18.65%  "
$ws.Range("E31").Value = "  -6.69%  "
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("D33").Value = "'27.42"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("D34").Value = "'0.112"
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "'0.993"
$ws.Range("E36").Value = "  +3.87%  "
$ws.Range("D37").Value = "'5.84"
$ws.Range("E37").Value = "  +8.19%  "
$ws.Range("D38").Value = "'2.20"
$ws.Range("E38").Value = "  +15.79%  "
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  +16.11%  "
$ws.Range("E42").Value = "  +7.95%  "
$ws.Range("D43").Value = "'43.61"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").Value = "'8.48"
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("D45").Value = "'398.84"
$ws.Range("E45").Value = "  +14.40%  "
$ws.Range("D46").Value = "'0.0364"
$ws.Range("E46").Value = "  +7.81%  "
$ws.Range("D47").Value = "2.807.18"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "'134.22"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D50").Value = "'24.05"
$ws.Range("E50").Value = "  +11.20%  "
$ws.Range("E51").Value = "  +4.37%  "
